# Append four new transaction rows (34-37) to the data table on the single
# worksheet, then leave the sheet scrolled/selected the way the author left
# it (top visible row ~7, with Q28:S40 selected) before saving.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data rows -------------------------------------------------------
# Row 34: Withdrawal / Credit Card / Tradeprof / 269.08
$ws.Range("E34").Value = "Withdrawal"
$ws.Range("N34").Value = "Credit Card"
$ws.Range("P34").Value = "Tradeprof"
$ws.Range("T34").Value = 269.08

# Row 35: Deposit / Crypto / ETH / 1578.4362000000001
$ws.Range("E35").Value = "Deposit"
$ws.Range("N35").Value = "Crypto"
$ws.Range("P35").Value = "ETH"
$ws.Range("T35").Value = 1578.4362000000001

# Row 36: Withdrawal / Credit Card / Tradeprof / 269.27499999999998
$ws.Range("E36").Value = "Withdrawal"
$ws.Range("N36").Value = "Credit Card"
$ws.Range("P36").Value = "Tradeprof"
$ws.Range("T36").Value = 269.27499999999998

# Row 37: Withdrawal / Credit Card / Sipay / 250
$ws.Range("E37").Value = "Withdrawal"
$ws.Range("N37").Value = "Credit Card"
$ws.Range("P37").Value = "Sipay"
$ws.Range("T37").Value = 250

# --- Restore the view the workbook was left in ---------------------------
# Scroll so row 7 is the first visible row (topLeftCell = A7).
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1

# Select the range the author had highlighted when the file was saved.
$ws.Range("Q28:S40").Select()
